$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1357
    $ws.Range("F4").Value = 7
    $ws.Range("F7").Value = 35
    $ws.Range("F8").Value = 182
}
